$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.356.15"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "1.596.70"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "211.84"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "19.13"
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("D11").Value = "0.0855"
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").Value = "1.820.45"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").Value = "1.591.37"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("D16").Value = "63.46"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").Value = "26.324.15"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "229.58"
$ws.Range("E18").Value = "  +7.30%  "
$ws.Range("D19").Value = "7.66"
$ws.Range("E19").Value = "  +4.13%  "
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("E23").Value = "  +2.70%  "
$ws.Range("D24").Value = "8.93"
$ws.Range("E24").Value = "  -1.29%  "
$ws.Range("D25").Value = "146.49"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "0.113"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").Value = "15.38"
$ws.Range("E29").Value = "  +1.88%  "
$ws.Range("D30").Value = "0.0495"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("D32").Value = "1.502.98"
$ws.Range("E32").Value = "  +5.49%  "
$ws.Range("E33").Value = "  +1.19%  "
$ws.Range("E34").Value = "  -0.84%  "
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").Value = "0.570"
$ws.Range("E37").Value = "  -3.09%  "
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("D39").Value = "0.818"
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("D40").Value = "5.78"
$ws.Range("E40").Value = "  -2.01%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  -5.00%  "
$ws.Range("D43").Value = "2.17"
$ws.Range("E43").Value = "  +1.96%  "
$ws.Range("D44").Value = "1.733.28"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("D46").Value = "60.66"
$ws.Range("D47").Value = "88.47"
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("D50").Value = "0.0958"
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("E51").Value = "  +0.11%  "
